# Generate Report for Handback
# - Marks the bfc561e4.../f5e0d5bd... rows as handed back (Status text + Latest
#   Handback DateTime) on both the Overview sheet and the per-locale sheets.
# - Populates the (previously empty) "Latest Target File" / "Latest Handback
#   File" columns (F/G) on the zh-cn and de-de sheets, with hyperlinks that
#   mirror the existing Source File Name / Latest Handoff File links.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: update the per-locale status cells (same shared text as the
# "Status" column on the locale sheets).
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $statusHandedBack
$ov.Range("C2").Value = $statusHandedBack
$ov.Range("B3").Value = $statusHandedBack
$ov.Range("C3").Value = $statusHandedBack

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Status column -> handed back
$zh.Range("C2").Value = $statusHandedBack
$zh.Range("C3").Value = $statusHandedBack

# Latest Handback DateTime
$zh.Range("H2").Value = "2016-03-20 14:55:17"
$zh.Range("H3").Value = "2016-03-20 14:55:17"

# Latest Target File (F) / Latest Handback File (G) for row 2 (bfc561e4...)
$zh.Range("F2").Value = "bfc561e4-cb0d-405b-ab5b-af1fba17e9ca.md"
$zh.Range("F2").Style = "HyperLink"
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/a0197ff2300a1c8db25873cb6738a9a78e1443a2/e2e/bfc561e4-cb0d-405b-ab5b-af1fba17e9ca.md", "", "", "bfc561e4-cb0d-405b-ab5b-af1fba17e9ca.md")

$zh.Range("G2").Value = "bfc561e4-cb0d-405b-ab5b-af1fba17e9ca.33bdfa43eab68736fa240edd557eb466bdf7bbf2.zh-cn.xlf"
$zh.Range("G2").Style = "HyperLink"
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3e43881c89f820da479fd557a9e928ba4479905c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bfc561e4-cb0d-405b-ab5b-af1fba17e9ca.33bdfa43eab68736fa240edd557eb466bdf7bbf2.zh-cn.xlf", "", "", "bfc561e4-cb0d-405b-ab5b-af1fba17e9ca.33bdfa43eab68736fa240edd557eb466bdf7bbf2.zh-cn.xlf")

# Latest Target File (F) / Latest Handback File (G) for row 3 (f5e0d5bd...)
$zh.Range("F3").Value = "f5e0d5bd-0f5b-4a88-990a-67a6ec7e5380.md"
$zh.Range("F3").Style = "HyperLink"
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/a0197ff2300a1c8db25873cb6738a9a78e1443a2/e2e/f5e0d5bd-0f5b-4a88-990a-67a6ec7e5380.md", "", "", "f5e0d5bd-0f5b-4a88-990a-67a6ec7e5380.md")

$zh.Range("G3").Value = "f5e0d5bd-0f5b-4a88-990a-67a6ec7e5380.c896566d183bcd0ee8f23f7838f5ad948e139868.zh-cn.xlf"
$zh.Range("G3").Style = "HyperLink"
$zh.Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3e43881c89f820da479fd557a9e928ba4479905c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f5e0d5bd-0f5b-4a88-990a-67a6ec7e5380.c896566d183bcd0ee8f23f7838f5ad948e139868.zh-cn.xlf", "", "", "f5e0d5bd-0f5b-4a88-990a-67a6ec7e5380.c896566d183bcd0ee8f23f7838f5ad948e139868.zh-cn.xlf")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Status column -> handed back
$de.Range("C2").Value = $statusHandedBack
$de.Range("C3").Value = $statusHandedBack

# Latest Handback DateTime
$de.Range("H2").Value = "2016-03-20 14:55:23"
$de.Range("H3").Value = "2016-03-20 14:55:23"

# Latest Target File (F) / Latest Handback File (G) for row 2 (bfc561e4...)
$de.Range("F2").Value = "bfc561e4-cb0d-405b-ab5b-af1fba17e9ca.md"
$de.Range("F2").Style = "HyperLink"
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/a0197ff2300a1c8db25873cb6738a9a78e1443a2/e2e/bfc561e4-cb0d-405b-ab5b-af1fba17e9ca.md", "", "", "bfc561e4-cb0d-405b-ab5b-af1fba17e9ca.md")

$de.Range("G2").Value = "bfc561e4-cb0d-405b-ab5b-af1fba17e9ca.33bdfa43eab68736fa240edd557eb466bdf7bbf2.de-de.xlf"
$de.Range("G2").Style = "HyperLink"
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/33a95d3635d847eefbd9a64e6794f2ad6d5fdc20/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bfc561e4-cb0d-405b-ab5b-af1fba17e9ca.33bdfa43eab68736fa240edd557eb466bdf7bbf2.de-de.xlf", "", "", "bfc561e4-cb0d-405b-ab5b-af1fba17e9ca.33bdfa43eab68736fa240edd557eb466bdf7bbf2.de-de.xlf")

# Latest Target File (F) / Latest Handback File (G) for row 3 (f5e0d5bd...)
$de.Range("F3").Value = "f5e0d5bd-0f5b-4a88-990a-67a6ec7e5380.md"
$de.Range("F3").Style = "HyperLink"
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/a0197ff2300a1c8db25873cb6738a9a78e1443a2/e2e/f5e0d5bd-0f5b-4a88-990a-67a6ec7e5380.md", "", "", "f5e0d5bd-0f5b-4a88-990a-67a6ec7e5380.md")

$de.Range("G3").Value = "f5e0d5bd-0f5b-4a88-990a-67a6ec7e5380.c896566d183bcd0ee8f23f7838f5ad948e139868.de-de.xlf"
$de.Range("G3").Style = "HyperLink"
$de.Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/33a95d3635d847eefbd9a64e6794f2ad6d5fdc20/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f5e0d5bd-0f5b-4a88-990a-67a6ec7e5380.c896566d183bcd0ee8f23f7838f5ad948e139868.de-de.xlf", "", "", "f5e0d5bd-0f5b-4a88-990a-67a6ec7e5380.c896566d183bcd0ee8f23f7838f5ad948e139868.de-de.xlf")
